# Update cryptos list values (price and 1h volume change) per row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.717.33"
$ws.Range("E2").Value = "  -0.60%  "

$ws.Range("D3").Value = "'2.920.25"
$ws.Range("E3").Value = "  +1.39%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'355.20"
$ws.Range("E5").Value = "  +0.93%  "

$ws.Range("D6").Value = "'110.56"
$ws.Range("E6").Value = "  -0.73%  "

$ws.Range("D7").Value = "'0.569"
$ws.Range("E7").Value = "  +2.00%  "

$ws.Range("D9").Value = "'0.630"
$ws.Range("E9").Value = "  +1.49%  "

$ws.Range("D10").Value = "'39.24"
$ws.Range("E10").Value = "  -1.97%  "

$ws.Range("D11").Value = "'0.0890"
$ws.Range("E11").Value = "  +3.15%  "

$ws.Range("E12").Value = "  +0.76%  "

$ws.Range("D13").Value = "'19.72"
$ws.Range("E13").Value = "  -1.54%  "

$ws.Range("D14").Value = "'7.91"
$ws.Range("E14").Value = "  +1.38%  "

$ws.Range("D15").Value = "'3.378.87"
$ws.Range("E15").Value = "  +1.27%  "

$ws.Range("D16").Value = "'2.909.28"
$ws.Range("E16").Value = "  +0.75%  "

$ws.Range("D17").Value = "'0.978"
$ws.Range("E17").Value = "  -1.50%  "

$ws.Range("D18").Value = "'51.760.03"
$ws.Range("E18").Value = "  -0.59%  "

$ws.Range("E19").Value = "  -2.30%  "

$ws.Range("D20").Value = "'3.27"
$ws.Range("E20").Value = "  -2.22%  "

$ws.Range("D21").Value = "'13.99"
$ws.Range("E21").Value = "  +0.45%  "

$ws.Range("D22").Value = "'0.0₃0983"
$ws.Range("E22").Value = "  +0.17%  "

$ws.Range("D23").Value = "'70.84"
$ws.Range("E23").Value = "  -0.19%  "

$ws.Range("D24").Value = "'270.34"
$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("E25").Value = "  +1.32%  "

$ws.Range("E26").Value = "  +12.15%  "

$ws.Range("D27").Value = "'27.09"
$ws.Range("E27").Value = "  +3.03%  "

$ws.Range("E28").Value = "  +0.09%  "

$ws.Range("D29").Value = "'7.43"
$ws.Range("E29").Value = "  +15.70%  "

$ws.Range("E30").Value = "  +13.78%  "

$ws.Range("D31").Value = "'10.58"
$ws.Range("E31").Value = "  +0.49%  "

$ws.Range("D32").Value = "'38.48"
$ws.Range("E32").Value = "  -0.37%  "

$ws.Range("D33").Value = "'6.06"
$ws.Range("E33").Value = "  +2.17%  "

$ws.Range("D34").Value = "'52.32"
$ws.Range("E34").Value = "  -1.70%  "

$ws.Range("E35").Value = "  -3.65%  "

$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("E37").Value = "  -15.99%  "

$ws.Range("D38").Value = "'3.24"
$ws.Range("E38").Value = "  -1.81%  "

$ws.Range("D39").Value = "'18.40"
$ws.Range("E39").Value = "  -0.94%  "

$ws.Range("E40").Value = "  -0.73%  "

$ws.Range("D41").Value = "'2.73"
$ws.Range("E41").Value = "  +3.89%  "

$ws.Range("E42").Value = "  +1.57%  "

$ws.Range("D43").Value = "'22.82"
$ws.Range("E43").Value = "  +1.44%  "

$ws.Range("D44").Value = "'118.80"
$ws.Range("E44").Value = "  -2.36%  "

$ws.Range("E45").Value = "  -2.42%  "

$ws.Range("E46").Value = "  +0.11%  "

$ws.Range("D47").Value = "'3.44"
$ws.Range("E47").Value = "  -3.94%  "

$ws.Range("D48").Value = "'2.129.47"
$ws.Range("E48").Value = "  -2.87%  "

$ws.Range("D49").Value = "'0.250"
$ws.Range("E49").Value = "  -7.62%  "

$ws.Range("D50").Value = "'0.0333"
$ws.Range("E50").Value = "  +4.81%  "

$ws.Range("D51").Value = "'9.10"
$ws.Range("E51").Value = "  -0.37%  "
